$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Formula = "=""66.778.95"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E2").Value = "  -2.08%  "

$c = $ws.Range("D3")
$c.Formula = "=""3.474.89"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E3").Value = "  -2.37%  "

$ws.Range("E4").Value = "  +0.06%  "

$c = $ws.Range("D5")
$c.Formula = "=""600.75"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E5").Value = "  -3.26%  "

$c = $ws.Range("D6")
$c.Formula = "=""146.79"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E6").Value = "  -5.34%  "

$c = $ws.Range("D7")
$c.Formula = "=""3.470.29"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E7").Value = "  -2.44%  "

$ws.Range("E8").Value = "  +0.00%  "

$c = $ws.Range("D9")
$c.Formula = "=""0.478"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E9").Value = "  -2.46%  "

$ws.Range("E10").Value = "  -3.28%  "

$c = $ws.Range("D11")
$c.Formula = "=""7.53"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E11").Value = "  +2.47%  "

$c = $ws.Range("D12")
$c.Formula = "=""0.420"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E12").Value = "  -4.05%  "

$ws.Range("E13").Value = "  -4.59%  "

$c = $ws.Range("D14")
$c.Formula = "=""4.063.72"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E14").Value = "  -2.31%  "

$c = $ws.Range("D15")
$c.Formula = "=""31.11"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E15").Value = "  -6.16%  "

$c = $ws.Range("D16")
$c.Formula = "=""3.470.03"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E16").Value = "  -2.50%  "

$c = $ws.Range("D17")
$c.Formula = "=""66.747.20"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E17").Value = "  -2.02%  "

$ws.Range("E18").Value = "  +0.12%  "

$c = $ws.Range("D19")
$c.Formula = "=""6.37"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E19").Value = "  -6.04%  "

$c = $ws.Range("D20")
$c.Formula = "=""15.24"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E20").Value = "  -4.67%  "

$ws.Range("E21").Value = "  +0.40%  "

$c = $ws.Range("D22")
$c.Formula = "=""433.08"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E22").Value = "  -4.86%  "

$c = $ws.Range("D23")
$c.Formula = "=""0.604"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E23").Value = "  -6.07%  "

$c = $ws.Range("D24")
$c.Formula = "=""79.29"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E24").Value = "  +1.04%  "

$c = $ws.Range("D25")
$c.Formula = "=""1.00"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E25").Value = "  -0.01%  "

$c = $ws.Range("D26")
$c.Formula = "=""3.613.85"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E26").Value = "  -2.29%  "

$ws.Range("E27").Value = "  -7.93%  "

$c = $ws.Range("D28")
$c.Formula = "=""9.73"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E28").Value = "  -7.66%  "

$ws.Range("E29").Value = "  -8.30%  "

$ws.Range("E30").Value = "  -3.56%  "

$c = $ws.Range("D31")
$c.Formula = "=""1.57"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E31").Value = "  -7.07%  "

$c = $ws.Range("D32")
$c.Formula = "=""0.167"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E32").Value = "  -2.69%  "

$c = $ws.Range("D33")
$c.Formula = "=""1.00"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E33").Value = "  +0.11%  "

$c = $ws.Range("D34")
$c.Formula = "=""25.24"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E34").Value = "  -3.30%  "

$c = $ws.Range("D35")
$c.Formula = "=""3.468.50"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E35").Value = "  -2.34%  "

$c = $ws.Range("D36")
$c.Formula = "=""5.88"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E36").Value = "  -8.23%  "

$ws.Range("E37").Value = "  -6.86%  "

$ws.Range("E38").Value = "  +0.02%  "

$c = $ws.Range("D39")
$c.Formula = "=""7.83"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E39").Value = "  -5.33%  "

$ws.Range("E40").Value = "  -0.04%  "

$c = $ws.Range("D41")
$c.Formula = "=""175.15"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E41").Value = "  -2.00%  "

$c = $ws.Range("D42")
$c.Formula = "=""0.0876"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E42").Value = "  -4.56%  "

$ws.Range("E43").Value = "  -11.23%  "

$c = $ws.Range("D44")
$c.Formula = "=""5.37"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E44").Value = "  -4.01%  "

$c = $ws.Range("D45")
$c.Formula = "=""0.890"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E45").Value = "  -0.87%  "

$c = $ws.Range("D46")
$c.Formula = "=""46.23"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E46").Value = "  -0.86%  "

$c = $ws.Range("D47")
$c.Formula = "=""28.67"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E47").Value = "  -7.90%  "

$ws.Range("E48").Value = "  -8.28%  "

$c = $ws.Range("D49")
$c.Formula = "=""7.39"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E49").Value = "  -5.12%  "

$ws.Range("E50").Value = "  -9.04%  "

$c = $ws.Range("D51")
$c.Formula = "=""0.971"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range("E51").Value = "  -5.01%  "
